$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Row 109 and row 112 exchange their match data (columns F, H, J..V).
#    Columns A-E, G, I stay put since they are identical between the two
#    rows already (Index, pais, torneio, temporada, data_partida, goals).
# ---------------------------------------------------------------------------

# Note: plain `.Value` cannot be read back into a variable in this
# environment (it yields the property descriptor string, not the data),
# so `.Value2` is used for every read here; writes still use `.Value`.

# --- capture current ("before edit") values of row 109 ---
$f109 = $ws.Range("F109").Value2
$h109 = $ws.Range("H109").Value2
$j109 = $ws.Range("J109").Value2
$k109 = $ws.Range("K109").Value2
$l109 = $ws.Range("L109").Value2
$m109 = $ws.Range("M109").Value2
$n109 = $ws.Range("N109").Value2
$o109 = $ws.Range("O109").Value2
$p109 = $ws.Range("P109").Value2
$q109 = $ws.Range("Q109").Value2
$r109 = $ws.Range("R109").Value2
$s109 = $ws.Range("S109").Value2
$t109 = $ws.Range("T109").Value2
$u109 = $ws.Range("U109").Value2
$v109 = $ws.Range("V109").Value2

# --- capture current ("before edit") values of row 112 ---
$f112 = $ws.Range("F112").Value2
$h112 = $ws.Range("H112").Value2
$j112 = $ws.Range("J112").Value2
$k112 = $ws.Range("K112").Value2
$l112 = $ws.Range("L112").Value2
$m112 = $ws.Range("M112").Value2
$n112 = $ws.Range("N112").Value2
$o112 = $ws.Range("O112").Value2
$p112 = $ws.Range("P112").Value2
$q112 = $ws.Range("Q112").Value2
$r112 = $ws.Range("R112").Value2
$s112 = $ws.Range("S112").Value2
$t112 = $ws.Range("T112").Value2
$u112 = $ws.Range("U112").Value2
$v112 = $ws.Range("V112").Value2

# --- write row 112's original values into row 109 ---
$ws.Range("F109").Value = $f112
$ws.Range("H109").Value = $h112
$ws.Range("J109").Value = $j112
$ws.Range("K109").Value = $k112
$ws.Range("L109").Value = $l112
$ws.Range("M109").Value = $m112
$ws.Range("N109").Value = $n112
$ws.Range("O109").Value = $o112
$ws.Range("P109").Value = $p112
$ws.Range("Q109").Value = $q112
$ws.Range("R109").Value = $r112
$ws.Range("S109").Value = $s112
$ws.Range("T109").Value = $t112
$ws.Range("U109").Value = $u112
$ws.Range("V109").Value = $v112

# --- write row 109's original values into row 112 ---
$ws.Range("F112").Value = $f109
$ws.Range("H112").Value = $h109
$ws.Range("J112").Value = $j109
$ws.Range("K112").Value = $k109
$ws.Range("L112").Value = $l109
$ws.Range("M112").Value = $m109
$ws.Range("N112").Value = $n109
$ws.Range("O112").Value = $o109
$ws.Range("P112").Value = $p109
$ws.Range("Q112").Value = $q109
$ws.Range("R112").Value = $r109
$ws.Range("S112").Value = $s109
$ws.Range("T112").Value = $t109
$ws.Range("U112").Value = $u109
$ws.Range("V112").Value = $v109

# ---------------------------------------------------------------------------
# 2) Append a new row 113 with a new match (Usti nad Labem vs Banik Most-Sous)
#    Copy formatting from row 112 first (borders/date styles etc.), then
#    overwrite with the new row's values.
# ---------------------------------------------------------------------------
$ws.Range("A112:V112").Copy($ws.Range("A113:V113"))

$ws.Range("A113").Value = 112
$ws.Range("B113").Value = "czech-republic"
$ws.Range("C113").Value = "cfl-group-b"
$ws.Range("D113").Value = "2023-2024"
$ws.Range("E113").Value = 45240.75
$ws.Range("F113").Value = "Usti nad Labem"
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = "Banik Most-Sous"
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1.81
$ws.Range("K113").Value = "10/11/2023 11:42"
$ws.Range("L113").Value = 1.43
$ws.Range("M113").Value = "10/11/2023 17:52"
$ws.Range("N113").Value = 3.66
$ws.Range("O113").Value = "10/11/2023 11:42"
$ws.Range("P113").Value = 4.77
$ws.Range("Q113").Value = "10/11/2023 17:52"
$ws.Range("R113").Value = 3.62
$ws.Range("S113").Value = "10/11/2023 11:42"
$ws.Range("T113").Value = 5.61
$ws.Range("U113").Value = "10/11/2023 17:52"
$ws.Range("V113").Value = "https://www.betexplorer.com/football/czech-republic/cfl-group-b/usti-nad-labem-banik-most-sous/xdnniW8C/"
